# Apply the commit: add a new customer row at the top of the data table
# (row 2), pushing all existing rows down by one, and correct the
# "Dư nợ" (outstanding debt) value for "Trương thị kiều" (now row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above current row 2 (the first data row),
# shifting all data below it down by one row.
$ws.Rows.Item(2).Insert()

# Fill in the new customer record in row 2.
$ws.Range("A2").Value = "KH"
$ws.Range("B2").Value = 440
$ws.Range("C2").Value = "Nguyễn Thị Như Ý"
$ws.Range("D2").Value = "CẦN THƠ"
# Phone numbers are stored as text (to keep the leading zero) -
# format the cell as Text before writing the digit string.
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "0966680674"
$ws.Range("I2").Value = 27000000
$ws.Range("J2").Value = 8000000

# Fix the outstanding debt ("Dư nợ") amount for Trương thị kiều, now on row 4.
$ws.Range("J4").Value = 26500000
